$d = $word.ActiveDocument

# Locate the point right after "On " (before "8.1.26") — this is where the
# _GoBack bookmark needs to move to. Using Find keeps this independent of
# hardcoded character offsets.
$anchor = $d.Content
$anchor.Find.Execute("On ") | Out-Null
$target = $d.Range($anchor.End, $anchor.End)

# Re-adding a bookmark with an existing name moves it (and splits the
# surrounding run so the bookmarkStart/bookmarkEnd pair lands at the new,
# now-collapsed, range) — exactly what's needed to relocate _GoBack from
# right after "classwork" to right after "On ".
$d.Bookmarks.Add("_GoBack", $target) | Out-Null
